$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.054482436495485
$ws.Range("D2").Value = 1.060446119281124
$ws.Range("E2").Value = 1.050990986245248
$ws.Range("F2").Value = 1.068265494641657
$ws.Range("I2").Value = 1.037987471969305
$ws.Range("J2").Value = 1.059494301194172
$ws.Range("K2").Value = 1.063172877107181
$ws.Range("L2").Value = 1.053743728498696
$ws.Range("M2").Value = 1.070971139364868
$ws.Range("N2").Value = 1.060998904627563

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.056030239185165
$ws.Range("D3").Value = 1.061888008690297
$ws.Range("E3").Value = 1.05232729638981
$ws.Range("F3").Value = 1.069876919874395
$ws.Range("I3").Value = 1.038270075816362
$ws.Range("J3").Value = 1.060690506652528
$ws.Range("K3").Value = 1.064427609904409
$ws.Range("L3").Value = 1.054891247837263
$ws.Range("M3").Value = 1.07239653348784
$ws.Range("N3").Value = 1.062196808834875

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.057030013025923
$ws.Range("D4").Value = 1.062819611122021
$ws.Range("E4").Value = 1.053190660904914
$ws.Range("F4").Value = 1.070918455706163
$ws.Range("I4").Value = 1.038450774612747
$ws.Range("J4").Value = 1.061462387340376
$ws.Range("K4").Value = 1.065237591656355
$ws.Range("L4").Value = 1.055631915650374
$ws.Range("M4").Value = 1.07331720832368
$ws.Range("N4").Value = 1.062969785681837

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.057449905854163
$ws.Range("D5").Value = 1.063210929917202
$ws.Range("E5").Value = 1.053553310574234
$ws.Range("F5").Value = 1.071356047129854
$ws.Range("I5").Value = 1.038526224014055
$ws.Range("J5").Value = 1.061786379264676
$ws.Range("K5").Value = 1.065577656534241
$ws.Range("L5").Value = 1.055942854701243
$ws.Range("M5").Value = 1.073703872558603
$ws.Range("N5").Value = 1.063294237711831

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.057520383796147
$ws.Range("D6").Value = 1.06327661511248
$ws.Range("E6").Value = 1.053614183134978
$ws.Range("F6").Value = 1.071429505052509
$ws.Range("I6").Value = 1.038538862074563
$ws.Range("J6").Value = 1.061840749411679
$ws.Range("K6").Value = 1.065634728647893
$ws.Range("L6").Value = 1.055995037246807
$ws.Range("M6").Value = 1.073768772685536
$ws.Range("N6").Value = 1.063348685070679

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.057035625262869
$ws.Range("D7").Value = 1.062824841218099
$ws.Range("E7").Value = 1.053195507852622
$ws.Range("F7").Value = 1.070924303875627
$ws.Range("I7").Value = 1.038451784798023
$ws.Range("J7").Value = 1.061466718519232
$ws.Range("K7").Value = 1.065242137385377
$ws.Range("L7").Value = 1.055636072145485
$ws.Range("M7").Value = 1.073322376463142
$ws.Range("N7").Value = 1.062974123011464

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.055005892963334
$ws.Range("D8").Value = 1.060933706191777
$ws.Range("E8").Value = 1.051442875323422
$ws.Range("F8").Value = 1.068810329459765
$ws.Range("I8").Value = 1.038083428634659
$ws.Range("J8").Value = 1.059899012402581
$ws.Range("K8").Value = 1.063597320416795
$ws.Range("L8").Value = 1.054131925606957
$ws.Range("M8").Value = 1.07145320434981
$ws.Range("N8").Value = 1.061404190572308

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.051415386574669
$ws.Range("D9").Value = 1.0575902547965
$ws.Range("E9").Value = 1.048344123443946
$ws.Range("F9").Value = 1.065075941863568
$ws.Range("I9").Value = 1.037417673258272
$ws.Range("J9").Value = 1.057119793812636
$ws.Range("K9").Value = 1.060683974375679
$ws.Range("L9").Value = 1.051466952024531
$ws.Range("M9").Value = 1.068146498928457
$ws.Range("N9").Value = 1.058621025173184

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04901184370997
$ws.Range("D10").Value = 1.055353400576256
$ws.Range("E10").Value = 1.046270883022984
$ws.Range("F10").Value = 1.062579568055084
$ws.Range("I10").Value = 1.03696251149623
$ws.Range("J10").Value = 1.055255320747745
$ws.Range("K10").Value = 1.058731247520011
$ws.Range("L10").Value = 1.049680181888776
$ws.Range("M10").Value = 1.065932813610963
$ws.Range("N10").Value = 1.056753904342683

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047968617636877
$ws.Range("D11").Value = 1.054382841897471
$ws.Range("E11").Value = 1.045371296120907
$ws.Range("F11").Value = 1.061496882717916
$ws.Range("I11").Value = 1.036762708072706
$ws.Range("J11").Value = 1.054445123585315
$ws.Range("K11").Value = 1.057883108768819
$ws.Range("L11").Value = 1.048904002881296
$ws.Range("M11").Value = 1.064971969224934
$ws.Range("N11").Value = 1.055942556607352

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.047580734243854
$ws.Range("D12").Value = 1.054022025633275
$ws.Range("E12").Value = 1.045036861813877
$ws.Range("F12").Value = 1.061094453551759
$ws.Range("I12").Value = 1.036688081814976
$ws.Range("J12").Value = 1.05414374120605
$ws.Range("K12").Value = 1.05756767394362
$ws.Range("L12").Value = 1.048615312877311
$ws.Range("M12").Value = 1.064614713246173
$ws.Range("N12").Value = 1.055640746230539

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.047663954043785
$ws.Range("D13").Value = 1.054099436012551
$ws.Range("E13").Value = 1.045108612258704
$ws.Range("F13").Value = 1.061180788432874
$ws.Range("I13").Value = 1.036704108009733
$ws.Range("J13").Value = 1.054208408751016
$ws.Range("K13").Value = 1.057635353945946
$ws.Range("L13").Value = 1.04867725529765
$ws.Range("M13").Value = 1.064691362110034
$ws.Range("N13").Value = 1.055705505610837

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047936562945808
$ws.Range("D14").Value = 1.054353023039951
$ws.Range("E14").Value = 1.045343657593081
$ws.Range("F14").Value = 1.061463623398821
$ws.Range("I14").Value = 1.036756547828446
$ws.Range("J14").Value = 1.054420220228555
$ws.Range("K14").Value = 1.05785704300266
$ws.Range("L14").Value = 1.048880147523869
$ws.Range("M14").Value = 1.064942445643138
$ws.Range("N14").Value = 1.055917617884969

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.048104475241886
$ws.Range("D15").Value = 1.054509225346801
$ws.Range("E15").Value = 1.04548843841397
$ws.Range("F15").Value = 1.061637851051644
$ws.Range("I15").Value = 1.036788803280803
$ws.Range("J15").Value = 1.054550665841983
$ws.Range("K15").Value = 1.057993579910866
$ws.Range("L15").Value = 1.049005105199028
$ws.Range("M15").Value = 1.065097099091554
$ws.Range("N15").Value = 1.056048248746131

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04908102610214
$ws.Range("D16").Value = 1.055417770734278
$ws.Range("E16").Value = 1.046330545752818
$ws.Range("F16").Value = 1.062651384794122
$ws.Range("I16").Value = 1.036975714375947
$ws.Range("J16").Value = 1.055309029741672
$ws.Range("K16").Value = 1.05878748030043
$ws.Range("L16").Value = 1.049731641103248
$ws.Range("M16").Value = 1.065996532341396
$ws.Range("N16").Value = 1.056807689609542

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.04969292014837
$ws.Range("D17").Value = 1.05598713907934
$ws.Range("E17").Value = 1.046858273654351
$ws.Range("F17").Value = 1.063286674866477
$ws.Range("I17").Value = 1.037092230086796
$ws.Range("J17").Value = 1.055783957769158
$ws.Range("K17").Value = 1.059284772313291
$ws.Range("L17").Value = 1.050186704351873
$ws.Range("M17").Value = 1.066560099563685
$ws.Range("N17").Value = 1.057283292089303

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.050049589412711
$ws.Range("D18").Value = 1.056319051002483
$ws.Range("E18").Value = 1.047165909495692
$ws.Range("F18").Value = 1.063657061733509
$ws.Range("I18").Value = 1.037159929899881
$ws.Range("J18").Value = 1.056060699087496
$ws.Range("K18").Value = 1.059574584512222
$ws.Range("L18").Value = 1.050451894759846
$ws.Range("M18").Value = 1.066888597125767
$ws.Range("N18").Value = 1.057560426412055

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.050171164339274
$ws.Range("D19").Value = 1.056432192270222
$ws.Range("E19").Value = 1.04727077526764
$ws.Range("F19").Value = 1.063783326055015
$ws.Range("I19").Value = 1.037182969448322
$ws.Range("J19").Value = 1.056155014129525
$ws.Range("K19").Value = 1.059673360863034
$ws.Range("L19").Value = 1.050542277307948
$ws.Range("M19").Value = 1.067000568991286
$ws.Range("N19").Value = 1.057654875392262

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.049627294422169
$ws.Range("D20").Value = 1.055926071066799
$ws.Range("E20").Value = 1.046801672010422
$ws.Range("F20").Value = 1.063218531611323
$ws.Range("I20").Value = 1.037079756153495
$ws.Range("J20").Value = 1.055733031110831
$ws.Range("K20").Value = 1.059231443487663
$ws.Range("L20").Value = 1.050137905280927
$ws.Range("M20").Value = 1.066499657130346
$ws.Range("N20").Value = 1.05723229310928

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047856297072276
$ws.Range("D21").Value = 1.054278356545043
$ws.Range("E21").Value = 1.045274450608111
$ws.Range("F21").Value = 1.061380343132656
$ws.Range("I21").Value = 1.036741116956603
$ws.Range("J21").Value = 1.054357859214561
$ws.Range("K21").Value = 1.057791772164909
$ws.Range("L21").Value = 1.048820411442168
$ws.Range("M21").Value = 1.064868517693012
$ws.Range("N21").Value = 1.055855168311183

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.046740580702345
$ws.Range("D22").Value = 1.05324058899513
$ws.Range("E22").Value = 1.044312557361141
$ws.Range("F22").Value = 1.060223024342215
$ws.Range("I22").Value = 1.036525825564399
$ws.Range("J22").Value = 1.053490689576189
$ws.Range("K22").Value = 1.056884284816989
$ws.Range("L22").Value = 1.047989833577582
$ws.Range("M22").Value = 1.063840893508426
$ws.Range("N22").Value = 1.054986767192458

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047332257112261
$ws.Range("D23").Value = 1.05379090139819
$ws.Range("E23").Value = 1.044822636090308
$ws.Range("F23").Value = 1.060836693867633
$ws.Range("I23").Value = 1.036640181586305
$ws.Range("J23").Value = 1.053950636549982
$ws.Range("K23").Value = 1.057365582718806
$ws.Range("L23").Value = 1.048430351338914
$ws.Range("M23").Value = 1.064385855164628
$ws.Range("N23").Value = 1.055447367343711

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.049656948599123
$ws.Range("D24").Value = 1.055953665663393
$ws.Range("E24").Value = 1.046827248409961
$ws.Range("F24").Value = 1.063249323131602
$ws.Range("I24").Value = 1.037085393396131
$ws.Range("J24").Value = 1.055756043529107
$ws.Range("K24").Value = 1.059255541260538
$ws.Range("L24").Value = 1.050159956223186
$ws.Range("M24").Value = 1.066526969149358
$ws.Range("N24").Value = 1.05725533820783

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052345317647873
$ws.Range("D25").Value = 1.058455972446507
$ws.Range("E25").Value = 1.049146499361079
$ws.Range("F25").Value = 1.066042526630272
$ws.Range("I25").Value = 1.037591774157936
$ws.Range("J25").Value = 1.057840312773198
$ws.Range("K25").Value = 1.061438962270022
$ws.Range("L25").Value = 1.052157666780288
$ws.Range("M25").Value = 1.069002950000282
$ws.Range("N25").Value = 1.059342567353315
